$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = '$ git commit -am "comments"'
$ws.Range("B13").Value = "directly commit changes from working area to repository. It skip staging area."

$ws.Range("A16").Select()
